$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 283 (pushes old rows 283..382 down to 284..383,
# Excel copies row formatting from the row above, matching dimension A1:T383)
$ws.Rows("283:283").Insert()

# Populate the newly inserted row 283 with the new weekly data record
$ws.Range("A283").Value = 7
$ws.Range("B283").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C283").Value = "Ñuble"
$ws.Range("D283").Value = 45215
$ws.Range("E283").Value = 16
$ws.Range("F283").Value = "Fruta"
$ws.Range("G283").Value = 100108
$ws.Range("H283").Value = "Tropicales y subtropicales"
$ws.Range("I283").Value = 100108005
$ws.Range("J283").Value = "Piña"
$ws.Range("K283").Value = "Caramelo"
$ws.Range("L283").Value = "Segunda"
$ws.Range("M283").Value = 50
$ws.Range("N283").Value = 23000
$ws.Range("O283").Value = 23000
$ws.Range("P283").Value = 23000
$ws.Range("Q283").Value = "`$/caja 14 unidades"
$ws.Range("R283").Value = "Ecuador"
$ws.Range("S283").Value = 1643
$ws.Range("T283").Value = 14
